# Kate.CheckList.MT.xlsx — "Revert "Revert "Completeion of id tag persName"""
#
# This re-applies a batch of ID-tag ("y" / "ND" / "Query" / "query" /
# "Query check" / "One of two names of a company") completions to column C
# (and a few column-D notes) of Sheet1, plus refreshes the sheet's view
# (zoom back to 100%, scroll position, active-cell selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell content updates (column C / D id-tag completions) -------------

# C3, C4 previously held "?" (shared string re-used from header row B2);
# together with all the other newly-tagged rows below they become "y".
$cells = @("C3","C4","C55","C58","C100","C108","C149","C150","C151","C152","C153","C154","C155","C156","C157","C158","C159","C161","C162","C163","C164","C165","C166","C168","C169","C170","C171","C173","C174","C175","C176","C177","C178","C179","C180","C181","C182","C183","C187","C188","C189","C190","C192","C193","C194","C195","C196","C197","C198","C199","C200","C201","C202","C204","C205","C206","C207","C208","C209","C210","C212","C305","C426","C431","C437","C445","C451","C458")
foreach ($cell in $cells) {
    $ws.Range($cell).Value = "y"
}

$cells = @("D101","D117","D124")
foreach ($cell in $cells) {
    $ws.Range($cell).Value = "Query"
}

$cells = @("C136","C137","C184","C185","C211","C246","C394","C461","C504")
foreach ($cell in $cells) {
    $ws.Range($cell).Value = "ND"
}

$cells = @("D167","D191")
foreach ($cell in $cells) {
    $ws.Range($cell).Value = "query"
}

$ws.Range("D300").Value = "Query check"
$ws.Range("D461").Value = "One of two names of a company"

# --- View refresh (best effort) ------------------------------------------
# Source workbook moved the scroll position/zoom/selection when these edits
# were made; reproduce what the COM object model exposes for that.
$ws.Range("C2").Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 2
